$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 115, shifting the
# existing "Zanahoria" weekly records (old rows 115-152) down to 117-154.
$ws.Rows.Item(115).Insert()
$ws.Rows.Item(115).Insert()

# New row 115: Primera, week of 2021-11-24 (serial 44524), Chillán
$ws.Range("A115").Value = 11
$ws.Range("B115").Value = "Vega Monumental Concepción"
$ws.Range("C115").Value = "Bíobío"
$ws.Range("D115").Value = 44524
$ws.Range("E115").Value = 8
$ws.Range("F115").Value = 100114013
$ws.Range("G115").Value = "Zanahoria"
$ws.Range("H115").Value = "Sin especificar"
$ws.Range("I115").Value = "Primera"
$ws.Range("J115").Value = 600
$ws.Range("K115").Value = 6500
$ws.Range("L115").Value = 7000
$ws.Range("M115").Value = 6750
$ws.Range("N115").Value = "`$/saco 20 kilos"
$ws.Range("O115").Value = "Chillán"
$ws.Range("P115").Value = 338
$ws.Range("Q115").Value = 20
$ws.Range("R115").Value = "Hortaliza"

# New row 116: Segunda, same week, Chillán
$ws.Range("A116").Value = 11
$ws.Range("B116").Value = "Vega Monumental Concepción"
$ws.Range("C116").Value = "Bíobío"
$ws.Range("D116").Value = 44524
$ws.Range("E116").Value = 8
$ws.Range("F116").Value = 100114013
$ws.Range("G116").Value = "Zanahoria"
$ws.Range("H116").Value = "Sin especificar"
$ws.Range("I116").Value = "Segunda"
$ws.Range("J116").Value = 600
$ws.Range("K116").Value = 6000
$ws.Range("L116").Value = 6000
$ws.Range("M116").Value = 6000
$ws.Range("N116").Value = "`$/saco 20 kilos"
$ws.Range("O116").Value = "Chillán"
$ws.Range("P116").Value = 300
$ws.Range("Q116").Value = 20
$ws.Range("R116").Value = "Hortaliza"
